# Commit: "Removing 'star' for wildcard"
# Plan: a star will indicate wildcards and a "fullword" marker will indicate
# an exact/whole-word match, but for now every list entry is just trimmed
# down to the start-of-word fragment that should match as a prefix.
#
# Net effect on the "Categories" sheet:
#   A2: "Happy"        -> "Hap"
#   A4: "Funny"        -> "Funn"
#   A5: "Intellingent" -> "Intelli"
# (B column values - Dog / Cat / Human - and the Sad entry are unaffected.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Hap"
$ws.Range("A4").Value = "Funn"
$ws.Range("A5").Value = "Intelli"

# Move the selection/active cell to A3 (previously it pointed past the used
# range, at A6).
[void]$ws.Range("A3").Select()
